$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("hpi")
$ws.Range("A2").Value = "Presence of skin changes such as sclerodactyly or telangiectasia"
$ws.Range("B2").Value = "These findings are characteristic of CREST syndrome and are not typically associated with food impaction."
$ws.Range("C2").Value = "Acute onset of severe chest pain or discomfort"
$ws.Range("D2").Value = "Acute chest pain is more indicative of food impaction than the gradual symptoms of CREST syndrome."
$ws.Range("A3").Value = "History of Raynaud's phenomenon"
$ws.Range("B3").Value = "Raynaud's phenomenon is a common feature of CREST syndrome and is not seen in food impaction."
$ws.Range("C3").Value = "Recent history of eating solid food that is difficult to swallow"
$ws.Range("D3").Value = "This history is directly related to food impaction and not typical for CREST syndrome."
$ws.Range("A4").Value = "Gradual onset of dysphagia (difficulty swallowing)"
$ws.Range("B4").Value = "Dysphagia in CREST syndrome is often progressive due to esophageal motility issues, unlike acute food impaction."
$ws.Range("C4").Value = "Presence of visible food bolus on examination"
$ws.Range("D4").Value = "A visible food bolus is a direct indicator of food impaction, not seen in CREST syndrome."
$ws.Range("A5").Value = "Associated symptoms of gastrointestinal reflux or heartburn"
$ws.Range("B5").Value = "These symptoms are more common in CREST syndrome due to esophageal involvement compared to food impaction."
$ws.Range("C5").Value = "History of previous episodes of food impaction"
$ws.Range("D5").Value = "Recurrent food impaction is a strong indicator of this condition, unlike CREST syndrome."
$ws.Range("A6").Value = "History of connective tissue disease"
$ws.Range("B6").Value = "A prior diagnosis of a connective tissue disease supports the diagnosis of CREST syndrome."
$ws.Range("C6").Value = "Immediate relief of symptoms after vomiting or regurgitation"
$ws.Range("D6").Value = "Relief after expulsion of food is characteristic of food impaction, not CREST syndrome."
$ws = $wb.Worksheets.Item("hist")
$ws.Range("C2").Value = "History of recent food bolus obstruction"
$ws.Range("D2").Value = "Recent episodes of food bolus obstruction are characteristic of food impaction."
$ws.Range("C3").Value = "History of dysphagia without systemic symptoms"
$ws.Range("D3").Value = "Dysphagia in food impaction typically occurs without the systemic symptoms seen in CREST syndrome."
$ws.Range("C4").Value = "No history of connective tissue disease"
$ws.Range("D4").Value = "Absence of connective tissue disease history supports a diagnosis of food impaction over CREST syndrome."
$ws.Range("C5").Value = "Previous endoscopic interventions for food impaction"
$ws.Range("D5").Value = "Prior endoscopic procedures for food impaction indicate a recurrent issue not related to CREST syndrome."
$ws.Range("A6").Value = "Use of immunosuppressive medications"
$ws.Range("B6").Value = "Patients with CREST syndrome may be on immunosuppressive therapy, which is not common in food impaction."
$ws.Range("C6").Value = "Use of medications for gastroesophageal reflux disease (GERD)"
$ws.Range("D6").Value = "Medications for GERD may indicate a history of esophageal issues more typical of food impaction."
$ws = $wb.Worksheets.Item("soc")
$ws.Range("D4").Value = "Prior episodes of food impaction suggest a higher likelihood of recurrence rather than CREST syndrome."
$ws.Range("D5").Value = "Absence of autoimmune diseases in the family reduces the likelihood of CREST syndrome."
$ws.Range("C6").Value = "Active participation in high-risk eating behaviors"
$ws.Range("D6").Value = "Engaging in behaviors that increase the risk of choking or impaction supports the diagnosis of food impaction."
$ws = $wb.Worksheets.Item("obj")
$ws.Range("D2").Value = "Acute chest pain is more characteristic of food impaction due to obstruction rather than the gradual symptoms of CREST syndrome."
$ws.Range("D3").Value = "A palpable food bolus is a direct indicator of food impaction, which is not seen in CREST syndrome."
$ws.Range("A4").Value = "Bilateral pulmonary fibrosis"
$ws.Range("B4").Value = "Pulmonary fibrosis can occur in CREST syndrome and is not a feature of food impaction."
$ws.Range("C4").Value = "Normal skin examination"
$ws.Range("D4").Value = "A normal skin exam suggests the absence of sclerodactyly or telangiectasia, favoring food impaction over CREST syndrome."
$ws.Range("B5").Value = "Dysphagia in CREST syndrome is due to esophageal motility issues, while in food impaction, it is typically acute and related to obstruction."
$ws.Range("C5").Value = "Localized tenderness in the chest or abdomen"
$ws.Range("D5").Value = "Localized tenderness is more indicative of food impaction rather than the systemic findings of CREST syndrome."
$ws.Range("A6").Value = "Decreased bowel sounds"
$ws.Range("B6").Value = "Decreased bowel sounds may indicate a motility disorder associated with CREST syndrome, unlike the normal or increased sounds in food impaction."
$ws.Range("C6").Value = "Normal respiratory exam"
$ws.Range("D6").Value = "A normal respiratory exam suggests no pulmonary complications, which would favor food impaction over CREST syndrome."
$ws = $wb.Worksheets.Item("test")
$ws.Range("D2").Value = "Direct observation of food bolus obstruction during endoscopy strongly indicates food impaction."
$ws.Range("D3").Value = "Normal motility suggests that the esophagus is functioning properly, which is inconsistent with achalasia."
$ws.Range("C4").Value = "Barium swallow study showing no significant esophageal dilation"
$ws.Range("D4").Value = "Absence of significant dilation suggests that the esophagus is not affected by achalasia, favoring food impaction."
$ws.Range("A5").Value = "Endoscopic findings of esophageal dilation"
$ws.Range("B5").Value = "Esophageal dilation is a common finding in achalasia, which is part of the CREST syndrome presentation."
$ws.Range("C5").Value = "Presence of acute symptoms after eating"
$ws.Range("D5").Value = "Acute symptoms following meals are more indicative of food impaction rather than chronic achalasia."
$ws.Range("A6").Value = "Imaging showing esophageal stasis"
$ws.Range("B6").Value = "Esophageal stasis is a result of impaired motility in achalasia, supporting the diagnosis of CREST syndrome."
$ws.Range("C6").Value = "Radiologic evidence of localized obstruction without dilation"
$ws.Range("D6").Value = "Localized obstruction without dilation is more consistent with food impaction than with achalasia."
